# Apply "Improving layout adding columns and result expander" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: only add underline style to A14 (content already correct) ---
$ws.Range("A14").Font.Underline = $true

# --- Row 15: finish the "turtle / side mount / sweep" entry (A15, already
#     correct text-wise, keep as-is; fill the rest of the row) ---
$ws.Cells.Item(15, 2).Value = "Feu"
$ws.Cells.Item(15, 3).Value = "00:01:30"
$ws.Cells.Item(15, 4).Value = "00:02:35"
$ws.Cells.Item(15, 5).Value = "Turtle"
$ws.Cells.Item(15, 6).Value = "Side mount"
$ws.Cells.Item(15, 7).Value = "sweep"
$ws.Cells.Item(15, 8).Value = "Arm lock variation."
$ws.Cells.Item(15, 9).Value = "Opponent is on the side of the turtle"
$ws.Cells.Item(15, 9).Style = "Normal"
$ws.Cells.Item(15, 10).Value = "Portuguese"

# --- Row 16: same video, next segment ---
$ws.Cells.Item(16, 1).Value = "https://www.youtube.com/watch?v=LjKhkKGF8Ug&ab_channel=FEUBJJ"
$ws.Cells.Item(16, 2).Value = "Feu"
$ws.Cells.Item(16, 3).Value = "00:02:54"
$ws.Cells.Item(16, 4).Value = "00:04:50"
$ws.Cells.Item(16, 5).Value = "Turtle"
$ws.Cells.Item(16, 6).Value = "Side mount"
$ws.Cells.Item(16, 7).Value = "sweep"
$ws.Cells.Item(16, 8).Value = "Arm lock variation."
$ws.Cells.Item(16, 9).Value = "Opponent is on front of the turtle"
$ws.Cells.Item(16, 10).Value = "Portuguese"

# --- Row 17: same video, third segment (no H/I info) ---
$ws.Cells.Item(17, 1).Value = "https://www.youtube.com/watch?v=LjKhkKGF8Ug&ab_channel=FEUBJJ"
$ws.Cells.Item(17, 2).Value = "Feu"
$ws.Cells.Item(17, 3).Value = "00:04:50"
$ws.Cells.Item(17, 4).Value = "00:06:30"
$ws.Cells.Item(17, 5).Value = "Turtle"
$ws.Cells.Item(17, 6).Value = "North South"
$ws.Cells.Item(17, 7).Value = "sweep"
$ws.Cells.Item(17, 10).Value = "Portuguese"
$ws.Cells.Item(17, 9).Font.Underline = $true

# --- Row 18: new video with hyperlink on the Link cell ---
$ws.Cells.Item(18, 1).Value = "https://www.youtube.com/watch?v=ojvH99btFYo&ab_channel=BernardoFariaBJJFanatics"
[void]$ws.Hyperlinks.Add($ws.Cells.Item(18, 1), "https://www.youtube.com/watch?v=ojvH99btFYo&ab_channel=BernardoFariaBJJFanatics")
$ws.Cells.Item(18, 2).Value = "Leonardo Nogueira"
$ws.Cells.Item(18, 3).Value = "00:01:22"
$ws.Cells.Item(18, 4).Value = "00:08:28"
$ws.Cells.Item(18, 5).Value = "Half Guard"
$ws.Cells.Item(18, 6).Value = "Guard"
$ws.Cells.Item(18, 7).Value = "Sweep"
$ws.Cells.Item(18, 9).Value = "Basic principles of the half gard, with knee shield"
$ws.Cells.Item(18, 10).Value = "English"

# --- Row 19: new video with hyperlink on the Link cell ---
$ws.Cells.Item(19, 1).Value = "https://www.youtube.com/watch?v=usFjw23WsMI&ab_channel=GRACIEMAG"
[void]$ws.Hyperlinks.Add($ws.Cells.Item(19, 1), "https://www.youtube.com/watch?v=usFjw23WsMI&ab_channel=GRACIEMAG")
$ws.Cells.Item(19, 2).Value = "Carlson Gracie Jr"
$ws.Cells.Item(19, 3).Value = "00:01:07"
$ws.Cells.Item(19, 4).Value = "00:02:19"
$ws.Cells.Item(19, 5).Value = "Spider Guard"
$ws.Cells.Item(19, 6).Value = "Side mount"
$ws.Cells.Item(19, 7).Value = "Passing"
$ws.Cells.Item(19, 8).Value = "Finalização Relogio"
$ws.Cells.Item(19, 10).Value = "Portuguese"

# --- Row 20: new video, result expander entry ---
$ws.Cells.Item(20, 1).Value = "https://www.youtube.com/watch?v=lYGLz0sagH8&ab_channel=GracieBarra"
$ws.Cells.Item(20, 2).Value = "Marcio Feitosa & Marcelo Souza"
$ws.Cells.Item(20, 3).Value = "00:05:00"
$ws.Cells.Item(20, 4).Value = "00:05:31"
$ws.Cells.Item(20, 7).Value = "Sparring"
$ws.Cells.Item(20, 9).Value = "Open Guard,"

# --- Update the active selection to match the final state ---
[void]$ws.Range("C6").Select()

Write-Host "Edit applied"
